$wb = $excel.ActiveWorkbook

# --- Rename the original sheet (Hoja1 -> Hoja2) and insert the two new sheets before it ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Name = "Hoja2"

$hoja2 = $wb.Worksheets.Item("Hoja2")
$imagesPool = $wb.Worksheets.Add($hoja2)
$imagesPool.Name = "ImagesPool"

$hoja2b = $wb.Worksheets.Item("Hoja2")
$images = $wb.Worksheets.Add($hoja2b)
$images.Name = "Images"

# --- Enter cell text in the same order the original author typed it, so that the ---
# --- shared-string table comes out in the same sequence as the authored workbook.  ---
$imagesPool.Range("A1:A3").NumberFormat = "@"

$imagesPool.Range("A1").Value() = "id"

$images.Range("A1").Value() = "ID"
$images.Range("B1").Value() = "PATH"
$images.Range("D1").Value() = "DESCRIPTION"
$images.Range("C1").Value() = "NAME"

$imagesPool.Range("B1").Value() = "path"
$imagesPool.Range("C1").Value() = "name"
$imagesPool.Range("D1").Value() = "description"

$imagesPool.Range("A2").Value() = "001"
$imagesPool.Range("A3").Value() = "td02"
$imagesPool.Range("B2").Value() = "Blegh"
$imagesPool.Range("B3").Value() = "Dani"
$imagesPool.Range("C2").Value() = "Ble"
$imagesPool.Range("C3").Value() = "Daniela"
$imagesPool.Range("D2").Value() = "This is a Bleg"
$imagesPool.Range("D3").Value() = "This is a Dani"

$tblPool = $imagesPool.ListObjects.Add(1, $imagesPool.Range("A1:D3"), $null, 1)
$tblPool.Name = "Tabla3"

# --- Fill the "Images" sheet body: mirrors Hoja2's Text/Number/int data, with a blank DESCRIPTION column ---
$imgData = @(
  @("WenaWena", 2, 4),
  @("WenaWena", 2, 5),
  @("WenaWena", 2, 10),
  @("WenaWena", 2, 123),
  @("WenaWena", 2, 46),
  @("Elio", 22, 13),
  @("Elio", 22, 180),
  @("Elio", 22, 123),
  @("Jaime", 10, 1),
  @("Jaime", 10, 3),
  @("Jaime", 10, 5),
  @("Feña", 1, 1)
)

$row = 2
foreach ($item in $imgData) {
  $images.Cells.Item($row, 1).Value() = $item[0]
  $images.Cells.Item($row, 2).Value() = $item[1]
  $images.Cells.Item($row, 3).Value() = $item[2]
  $images.Cells.Item($row, 4).Value() = ""
  $row = $row + 1
}

$tblImages = $images.ListObjects.Add(1, $images.Range("A1:D13"), $null, 1)
$tblImages.Name = "Tabla13"

# --- Fix up selections / active tab ---
$hoja2c = $wb.Worksheets.Item("Hoja2")
$hoja2c.Select()
$hoja2c.Range("A2:C13").Select()

$imagesPool.Select()
